$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths for columns C (TITLE) and D (COUNTRY)
$ws.Columns.Item(3).ColumnWidth = 60.17
$ws.Columns.Item(4).ColumnWidth = 88.17

# The yellow PREMIUM highlight currently on E2 moves down to the new row 25 below.
# Copy E2's current format (while it still has the highlight) before we overwrite it.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Now clear E2's highlight - row 2 is no longer a premium opportunity
$ws.Range("E2").ClearFormats()

# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = '1331613'
$ws.Range("B2").Value = 'https://aiesec.org/opportunity/global-talent/1331613'
$ws.Range("C2").Value = 'Social Media & Content Creator'
$ws.Range("D2").Value = 'Giza, El Omraniya, Giza Governorate, Egypt'
$ws.Range("E2").Value = 'No'
$ws.Range("F2").Value = '0 applicants'
$ws.Range("G2").Value = '9 - 12 Weeks'
$ws.Range("H2").Value = 'EG scout shop'

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = '1331612'
$ws.Range("B3").Value = 'https://aiesec.org/opportunity/global-talent/1331612'
$ws.Range("C3").Value = 'Web Developer'
$ws.Range("D3").Value = 'Giza, El Omraniya, Giza Governorate, Egypt'
$ws.Range("E3").Value = 'No'
$ws.Range("F3").Value = '0 applicants'
$ws.Range("G3").Value = '9 - 12 Weeks'
$ws.Range("H3").Value = 'EG scout shop'

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = '1331309'
$ws.Range("B4").Value = 'https://aiesec.org/opportunity/global-talent/1331309'
$ws.Range("C4").Value = 'Interior designer'
$ws.Range("D4").Value = 'Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt'
$ws.Range("E4").Value = 'No'
$ws.Range("F4").Value = '0 applicants'
$ws.Range("G4").Value = '9 - 12 Weeks'
$ws.Range("H4").Value = 'AI design'

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = '1331158'
$ws.Range("B5").Value = 'https://aiesec.org/opportunity/global-talent/1331158'
$ws.Range("C5").Value = 'Business Development Executive'
$ws.Range("D5").Value = 'Cairo, Cairo Governorate, Egypt'
$ws.Range("E5").Value = 'No'
$ws.Range("F5").Value = '3 applicants'
$ws.Range("G5").Value = '9 - 12 Weeks'
$ws.Range("H5").Value = 'Silverkey Technologies Egypt'

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = '1331148'
$ws.Range("B6").Value = 'https://aiesec.org/opportunity/global-talent/1331148'
$ws.Range("C6").Value = 'Digital Marketing Executive'
$ws.Range("D6").Value = 'Cairo, Cairo Governorate, Egypt'
$ws.Range("E6").Value = 'No'
$ws.Range("F6").Value = '5 applicants'
$ws.Range("G6").Value = '9 - 12 Weeks'
$ws.Range("H6").Value = 'Silverkey Technologies Egypt'

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = '1331071'
$ws.Range("B7").Value = 'https://aiesec.org/opportunity/global-talent/1331071'
$ws.Range("C7").Value = 'Growth Analytics & Funnel Optimization'
$ws.Range("D7").Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E7").Value = 'No'
$ws.Range("F7").Value = '1 applicant'
$ws.Range("G7").Value = '9 - 12 Weeks'
$ws.Range("H7").Value = 'Madaar'

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = '1331070'
$ws.Range("B8").Value = 'https://aiesec.org/opportunity/global-talent/1331070'
$ws.Range("C8").Value = 'Performance Marketing'
$ws.Range("D8").Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E8").Value = 'No'
$ws.Range("F8").Value = '3 applicants'
$ws.Range("G8").Value = '9 - 12 Weeks'
$ws.Range("H8").Value = 'Madaar'

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = '1331068'
$ws.Range("B9").Value = 'https://aiesec.org/opportunity/global-talent/1331068'
$ws.Range("C9").Value = 'Product Marketing'
$ws.Range("D9").Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E9").Value = 'No'
$ws.Range("F9").Value = '1 applicant'
$ws.Range("G9").Value = '9 - 12 Weeks'
$ws.Range("H9").Value = 'Madaar'

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = '1331067'
$ws.Range("B10").Value = 'https://aiesec.org/opportunity/global-talent/1331067'
$ws.Range("C10").Value = 'Business Development'
$ws.Range("D10").Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E10").Value = 'No'
$ws.Range("F10").Value = '1 applicant'
$ws.Range("G10").Value = '9 - 12 Weeks'
$ws.Range("H10").Value = 'Madaar'

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = '1331066'
$ws.Range("B11").Value = 'https://aiesec.org/opportunity/global-talent/1331066'
$ws.Range("C11").Value = 'Sales'
$ws.Range("D11").Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E11").Value = 'No'
$ws.Range("F11").Value = '0 applicants'
$ws.Range("G11").Value = '9 - 12 Weeks'
$ws.Range("H11").Value = 'Madaar'

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = '1331065'
$ws.Range("B12").Value = 'https://aiesec.org/opportunity/global-talent/1331065'
$ws.Range("C12").Value = 'Content & Brand Marketing'
$ws.Range("D12").Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E12").Value = 'No'
$ws.Range("F12").Value = '2 applicants'
$ws.Range("G12").Value = '9 - 12 Weeks'
$ws.Range("H12").Value = 'Madaar'

# Row 13
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = '1330036'
$ws.Range("B13").Value = 'https://aiesec.org/opportunity/global-talent/1330036'
$ws.Range("C13").Value = '[Accelerate Serbia] Structural engineer'
$ws.Range("D13").Value = 'Belgrade, Serbia'
$ws.Range("E13").Value = 'No'
$ws.Range("F13").Value = '25 applicants'
$ws.Range("G13").Value = '9 - 12 Weeks'
$ws.Range("H13").Value = 'Welt Inzenjering'

# Row 14
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = '1329104'
$ws.Range("B14").Value = 'https://aiesec.org/opportunity/global-talent/1329104'
$ws.Range("C14").Value = 'Sales Intern'
$ws.Range("D14").Value = 'Pannipitiya, Sri Lanka'
$ws.Range("E14").Value = 'No'
$ws.Range("F14").Value = '22 applicants'
$ws.Range("G14").Value = '9 - 12 Weeks'
$ws.Range("H14").Value = 'Frella International'

# Row 15
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = '1328557'
$ws.Range("B15").Value = 'https://aiesec.org/opportunity/global-talent/1328557'
$ws.Range("C15").Value = 'Arduino Developer'
$ws.Range("D15").Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E15").Value = 'No'
$ws.Range("F15").Value = '9 applicants'
$ws.Range("G15").Value = '9 - 12 Weeks'
$ws.Range("H15").Value = 'Techno square'

# Row 16
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = '1328553'
$ws.Range("B16").Value = 'https://aiesec.org/opportunity/global-talent/1328553'
$ws.Range("C16").Value = 'Business developer'
$ws.Range("D16").Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E16").Value = 'No'
$ws.Range("F16").Value = '22 applicants'
$ws.Range("G16").Value = '9 - 12 Weeks'
$ws.Range("H16").Value = 'I.C.Robotics'

# Row 17
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = '1328552'
$ws.Range("B17").Value = 'https://aiesec.org/opportunity/global-talent/1328552'
$ws.Range("C17").Value = 'Mobile applicatio'
$ws.Range("D17").Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E17").Value = 'No'
$ws.Range("F17").Value = '6 applicants'
$ws.Range("G17").Value = '9 - 12 Weeks'
$ws.Range("H17").Value = 'TAR - Company'

# Row 18
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = '1328548'
$ws.Range("B18").Value = 'https://aiesec.org/opportunity/global-talent/1328548'
$ws.Range("C18").Value = 'Sales'
$ws.Range("D18").Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E18").Value = 'No'
$ws.Range("F18").Value = '25 applicants'
$ws.Range("G18").Value = '9 - 12 Weeks'
$ws.Range("H18").Value = 'TAR - Company'

# Row 19
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = '1328547'
$ws.Range("B19").Value = 'https://aiesec.org/opportunity/global-talent/1328547'
$ws.Range("C19").Value = 'SEO'
$ws.Range("D19").Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E19").Value = 'No'
$ws.Range("F19").Value = '14 applicants'
$ws.Range("G19").Value = '9 - 12 Weeks'
$ws.Range("H19").Value = 'TAR - Company'

# Row 20
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = '1328545'
$ws.Range("B20").Value = 'https://aiesec.org/opportunity/global-talent/1328545'
$ws.Range("C20").Value = 'Graphic designer'
$ws.Range("D20").Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E20").Value = 'No'
$ws.Range("F20").Value = '14 applicants'
$ws.Range("G20").Value = '9 - 12 Weeks'
$ws.Range("H20").Value = 'TAR - Company'

# Row 21
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = '1328543'
$ws.Range("B21").Value = 'https://aiesec.org/opportunity/global-talent/1328543'
$ws.Range("C21").Value = 'UI/UX design'
$ws.Range("D21").Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E21").Value = 'No'
$ws.Range("F21").Value = '77 applicants'
$ws.Range("G21").Value = '9 - 12 Weeks'
$ws.Range("H21").Value = 'TAR - Company'

# Row 22
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = '1328365'
$ws.Range("B22").Value = 'https://aiesec.org/opportunity/global-talent/1328365'
$ws.Range("C22").Value = 'Graphic Designer'
$ws.Range("D22").Value = 'Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt'
$ws.Range("E22").Value = 'No'
$ws.Range("F22").Value = '8 applicants'
$ws.Range("G22").Value = '9 - 12 Weeks'
$ws.Range("H22").Value = 'ASG Engineering'

# Row 23
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = '1328363'
$ws.Range("B23").Value = 'https://aiesec.org/opportunity/global-talent/1328363'
$ws.Range("C23").Value = 'Marketing Specialist'
$ws.Range("D23").Value = 'Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt'
$ws.Range("E23").Value = 'No'
$ws.Range("F23").Value = '19 applicants'
$ws.Range("G23").Value = '9 - 12 Weeks'
$ws.Range("H23").Value = 'ASG Engineering'

# Row 24
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = '1328300'
$ws.Range("B24").Value = 'https://aiesec.org/opportunity/global-talent/1328300'
$ws.Range("C24").Value = 'content creator'
$ws.Range("D24").Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Range("E24").Value = 'No'
$ws.Range("F24").Value = '16 applicants'
$ws.Range("G24").Value = '9 - 12 Weeks'
$ws.Range("H24").Value = 'Markit'

# Row 25
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = '1327954'
$ws.Range("B25").Value = 'https://aiesec.org/opportunity/global-talent/1327954'
$ws.Range("C25").Value = 'Taste Hungary | [Only EU] Jr. Commercial Analyst (Finance)'
$ws.Range("D25").Value = 'Budapeste, Hungria'
$ws.Range("E25").Value = 'Yes'
$ws.Range("F25").Value = '21 applicants'
$ws.Range("G25").Value = '6 - 18 Months'
$ws.Range("H25").Value = 'EATON'

# Row 26
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = '1326381'
$ws.Range("B26").Value = 'https://aiesec.org/opportunity/global-talent/1326381'
$ws.Range("C26").Value = 'Business Developer'
$ws.Range("D26").Value = 'Sheraton Al Matar, El Nozha, Cairo Governorate, Egypt'
$ws.Range("E26").Value = 'No'
$ws.Range("F26").Value = '12 applicants'
$ws.Range("G26").Value = '3 - 6 Months'
$ws.Range("H26").Value = '12 applicants'

# Row 27
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = '1324596'
$ws.Range("B27").Value = 'https://aiesec.org/opportunity/global-talent/1324596'
$ws.Range("C27").Value = 'ACCOUNTANT'
$ws.Range("D27").Value = 'New Damietta City, Damietta El-Gadeeda City, New Damietta, Damietta Governorate, Egypt'
$ws.Range("E27").Value = 'No'
$ws.Range("F27").Value = '49 applicants'
$ws.Range("G27").Value = '3 - 6 Months'
$ws.Range("H27").Value = 'Business Haven Consultancy'

# Row 28
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = '1324593'
$ws.Range("B28").Value = 'https://aiesec.org/opportunity/global-talent/1324593'
$ws.Range("C28").Value = 'Graphic desgin'
$ws.Range("D28").Value = 'New Damietta City, Damietta El-Gadeeda City, New Damietta, Damietta Governorate, Egypt'
$ws.Range("E28").Value = 'No'
$ws.Range("F28").Value = '10 applicants'
$ws.Range("G28").Value = '3 - 6 Months'
$ws.Range("H28").Value = 'Business Haven Consultancy'

# Row 29
$ws.Range("A29").NumberFormat = "@"
$ws.Range("A29").Value = '1324592'
$ws.Range("B29").Value = 'https://aiesec.org/opportunity/global-talent/1324592'
$ws.Range("C29").Value = 'Digital marketing'
$ws.Range("D29").Value = 'New Damietta City, Damietta El-Gadeeda City, New Damietta, Damietta Governorate, Egypt'
$ws.Range("E29").Value = 'No'
$ws.Range("F29").Value = '29 applicants'
$ws.Range("G29").Value = '3 - 6 Months'
$ws.Range("H29").Value = 'Business Haven Consultancy'

# Row 30
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = '1322224'
$ws.Range("B30").Value = 'https://aiesec.org/opportunity/global-talent/1322224'
$ws.Range("C30").Value = 'Technical Support Sales Specialist'
$ws.Range("D30").Value = 'Pendik, Kaynarca, 34890 Pendik/İstanbul, Türkiye'
$ws.Range("E30").Value = 'No'
$ws.Range("F30").Value = '75 applicants'
$ws.Range("G30").Value = '6 - 18 Months'
$ws.Range("H30").Value = 'AYBEY ELEKTRONİK SANAYİ VE TİCARET ANONİM ŞİRKETİ'

# Row 31
$ws.Range("A31").NumberFormat = "@"
$ws.Range("A31").Value = '1320725'
$ws.Range("B31").Value = 'https://aiesec.org/opportunity/global-talent/1320725'
$ws.Range("C31").Value = 'International Educational Consultant'
$ws.Range("D31").Value = 'İstanbul, Türkiye'
$ws.Range("E31").Value = 'No'
$ws.Range("F31").Value = '57 applicants'
$ws.Range("G31").Value = '6 - 18 Months'
$ws.Range("H31").Value = 'JOHN AND JOHN EĞİTİM TEKNOLOJİ VE İNTERNET YATIRIMLARI LİMİT'

# Row 32
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = '1303857'
$ws.Range("B32").Value = 'https://aiesec.org/opportunity/global-talent/1303857'
$ws.Range("C32").Value = 'MARKETING & SALES / ABC TEKSTİL'
$ws.Range("D32").Value = 'Pınarkent, 20180 Pamukkale/Denizli, Türkiye'
$ws.Range("E32").Value = 'No'
$ws.Range("F32").Value = '160 applicants'
$ws.Range("G32").Value = '6 - 18 Months'
$ws.Range("H32").Value = 'ABC TEKSTİL SANAYİ VE TİCARET ANONİM ŞİRKETİ'

# Row 33
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = '1288588'
$ws.Range("B33").Value = 'https://aiesec.org/opportunity/global-talent/1288588'
$ws.Range("C33").Value = 'Technical Marketing Intern'
$ws.Range("D33").Value = 'Colombo, Sri Lanka'
$ws.Range("E33").Value = 'No'
$ws.Range("F33").Value = '47 applicants'
$ws.Range("G33").Value = '6 - 18 Months'
$ws.Range("H33").Value = 'Epigro Pvt Ltd'

